# Update hotel reviews data: fill in English_Reviews_num and Local_Rank
# for the hotel row on the "hotel_info" sheet.
#
# The target cells (G2, H2) currently hold empty-string values typed as
# text (t="s"). We want to replace them with the text values "1" and
# "28" respectively, keeping them as text (not numbers), matching the
# existing column's string typing, then drop back to the default
# "Normal" style so no extra formatting is left behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "28"
$ws.Range("H2").Style = "Normal"
